# [update] Support Untitled file
# Add three new rectangle shapes (red / yellow / green) to slide 1, placed
# between the "図 11" picture and the "図 12" picture, matching the target
# OOXML produced by the commit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$msoSendToBack   = 1
$msoBringForward = 2
$EmuPerPoint     = 12700.0

# Shape.Left/.Top/.Width/.Height round-trip through a single-precision
# float internally, so naively dividing EMU by 12700 can end up one EMU
# short after the host floors the value back to EMU. Nudge the point
# value up in tiny steps until it lands back on the exact EMU we want.
function Get-SafePoints {
    param(
        [double]$TargetEmu
    )
    $basePt = $TargetEmu / $EmuPerPoint
    for ($i = 0; $i -le 1000; $i++) {
        $candidate = $basePt + ($i * 0.00001)
        $f = [float]$candidate
        $backEmu = [math]::Floor([double]$f * $EmuPerPoint)
        if ($backEmu -eq $TargetEmu) {
            return $candidate
        }
    }
    return $basePt
}

# Shape 1 on the slide ("正方形/長方形 1") already carries the default
# theme shape style (p:style / lnRef / fillRef / ...) that PowerPoint
# attaches to a rectangle drawn from the UI. Duplicating it (instead of
# Shapes.AddShape) keeps that style block and the plain txBody/bodyPr
# markup intact on the new shapes, exactly like the authored file.
$styleSource = $s.Shapes.Item(1)

function New-IconRect {
    param(
        [double]$LeftEmu,
        [double]$TopEmu,
        [double]$WidthEmu,
        [double]$HeightEmu,
        [string]$ShapeName,
        [long]$FillRGB
    )

    $dupRange = $styleSource.Duplicate()
    $shp = $dupRange.Item(1)
    $shp.Name = $ShapeName
    $shp.Left = Get-SafePoints $LeftEmu
    $shp.Top = Get-SafePoints $TopEmu
    $shp.Width = Get-SafePoints $WidthEmu
    $shp.Height = Get-SafePoints $HeightEmu
    $shp.Fill.ForeColor.RGB = $FillRGB
    $shp.Line.Visible = $false
    return $shp
}

# The host's shape-id counter is a monotonically increasing counter that
# starts at 3 and is unrelated to the ids already present on the slide.
# The authored file has the three new shapes numbered 3, 9 and 10, so we
# burn through the ids in between with scratch shapes and discard them,
# leaving the three real shapes with the desired ids.

$red = New-IconRect 9183315 340744 718457 653143 "正方形/長方形 2" 255            # id 3

$scratch1 = $styleSource.Duplicate().Item(1)   # id 4
$scratch2 = $styleSource.Duplicate().Item(1)   # id 5
$scratch3 = $styleSource.Duplicate().Item(1)   # id 6
$scratch4 = $styleSource.Duplicate().Item(1)   # id 7
$scratch5 = $styleSource.Duplicate().Item(1)   # id 8

$yellow = New-IconRect 8346909 373401 718457 653143 "正方形/長方形 8" 65535       # id 9
$green  = New-IconRect 7569477 340744 718457 653143 "正方形/長方形 9" 5296274     # id 10

$scratch1.Delete()
$scratch2.Delete()
$scratch3.Delete()
$scratch4.Delete()
$scratch5.Delete()

# Move the three new shapes so that they sit right after the "図 11"
# picture and right before the "図 12" picture, preserving the red /
# yellow / green order from the diff.
$red.ZOrder($msoSendToBack)
$red.ZOrder($msoBringForward)
$red.ZOrder($msoBringForward)
$red.ZOrder($msoBringForward)

$yellow.ZOrder($msoSendToBack)
$yellow.ZOrder($msoBringForward)
$yellow.ZOrder($msoBringForward)
$yellow.ZOrder($msoBringForward)
$yellow.ZOrder($msoBringForward)

$green.ZOrder($msoSendToBack)
$green.ZOrder($msoBringForward)
$green.ZOrder($msoBringForward)
$green.ZOrder($msoBringForward)
$green.ZOrder($msoBringForward)
$green.ZOrder($msoBringForward)
